$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> IonQ, Inc. (IONQ)
$ws.Range("B2").Value = "IonQ, Inc."
$ws.Range("C2").Value = "IONQ"
$ws.Range("D2").Value = 54.26
$ws.Range("E2").Value = 64.7
$ws.Range("F2").Value = 15.7
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 60
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 63
$ws.Range("K2").Value = 57.7
$ws.Range("N2").Value = 52.47848103381103

# Row 3 -> D-Wave Quantum Inc. (QBTS)
$ws.Range("B3").Value = "D-Wave Quantum Inc."
$ws.Range("C3").Value = "QBTS"
$ws.Range("D3").Value = 27.8
$ws.Range("E3").Value = 63.8
$ws.Range("F3").Value = 24.06
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 66
$ws.Range("I3").Value = 66
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 57.1
$ws.Range("N3").Value = 52.47848103381103

# Row 4 -> Rigetti Computing, Inc. (RGTI) - values updated
$ws.Range("D4").Value = 28.66
$ws.Range("E4").Value = 60.1
$ws.Range("F4").Value = 12.06
$ws.Range("N4").Value = 52.47848103381103

# Row 5 -> International Business Machines (IBM) - values updated
$ws.Range("D5").Value = 308.03
$ws.Range("E5").Value = 53
$ws.Range("F5").Value = 1.59
$ws.Range("N5").Value = 52.47848103381103
